# Restore C10 ("R30" rule row, "From" column) on the Rules sheet back to 1
# (matches the target revision's stored value, replacing the prior 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
